# Generate Report for Handoff
# Adds a new row (row 3) describing the 77872216-... markdown file to the
# Overview, zh-cn and de-de sheets, mirroring the existing row 2 pattern.

$wb = $excel.ActiveWorkbook

$mdName   = '77872216-f9d0-429b-b33c-da609c8242d5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$mdPath   = 'e2e\77872216-f9d0-429b-b33c-da609c8242d5ooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooooo.md'
$status   = 'Ready for handoff'
$date2820 = '2016-08-28 02:28:20'
$zhcnXlf  = '77872216-f9d0-429b-b33c-da609c8242d5oooooooooooooooooooooooooooooooooooooooo.f46ca0f2172d523803582f58daa28b2783955079.zh-cn.xlf'
$date2816 = '2016-08-28 02:28:16'
$dedeXlf  = '77872216-f9d0-429b-b33c-da609c8242d5oooooooooooooooooooooooooooooooooooooooo.f46ca0f2172d523803582f58daa28b2783955079.de-de.xlf'

$hyperlinkTarget = 'https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4eaa172c15aad46d038c7f0eea0a61c8dd27f701/e2e/' + $mdName
$dateFormat = 'yyyy-mm-dd HH:mm:ss'

# ---------------------------------------------------------------------
# Sheet "Overview" (sheet1 / table3) -> new row A3:G3
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("A3").Value = $mdName
$wsOverview.Range("B3").Value = $mdPath
$wsOverview.Range("C3").Value = ".md"
$wsOverview.Range("D3").Value = ""
$wsOverview.Range("E3").Value = $status
$wsOverview.Range("F3").Value = $status
$wsOverview.Range("G3").Value = $date2820
$wsOverview.Range("G3").NumberFormat = $dateFormat

$wsOverview.Hyperlinks.Add($wsOverview.Range("B3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $mdPath)

$loOverview = $wsOverview.ListObjects.Item(1)
$loOverview.Resize($wsOverview.Range("A1:G3"))

# ---------------------------------------------------------------------
# Sheet "zh-cn" (sheet2 / table1) -> new row A3:P3
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("A3").Value = $mdName
$wsZhCn.Range("B3").Value = ".md"
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("D3").Value = "e2e"
$wsZhCn.Range("E3").Value = "ht"
$wsZhCn.Range("F3").Value = "False"
$wsZhCn.Range("G3").Value = $zhcnXlf
$wsZhCn.Range("H3").Value = $date2816
$wsZhCn.Range("H3").NumberFormat = $dateFormat
$wsZhCn.Range("I3").Value = ""
$wsZhCn.Range("J3").Value = ""
$wsZhCn.Range("K3").Value = "0001-01-01 00:00:00"
$wsZhCn.Range("K3").NumberFormat = $dateFormat
$wsZhCn.Range("L3").Value = ""
$wsZhCn.Range("M3").Value = "True"
$wsZhCn.Range("N3").Value = ""
$wsZhCn.Range("O3").Value = "False"
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("A3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $mdName)

$loZhCn = $wsZhCn.ListObjects.Item(1)
$loZhCn.Resize($wsZhCn.Range("A1:P3"))

# ---------------------------------------------------------------------
# Sheet "de-de" (sheet3 / table2) -> new row A3:P3
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("A3").Value = $mdName
$wsDeDe.Range("B3").Value = ".md"
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("D3").Value = "e2e"
$wsDeDe.Range("E3").Value = "ht"
$wsDeDe.Range("F3").Value = "False"
$wsDeDe.Range("G3").Value = $dedeXlf
$wsDeDe.Range("H3").Value = $date2820
$wsDeDe.Range("H3").NumberFormat = $dateFormat
$wsDeDe.Range("I3").Value = ""
$wsDeDe.Range("J3").Value = ""
$wsDeDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDeDe.Range("K3").NumberFormat = $dateFormat
$wsDeDe.Range("L3").Value = ""
$wsDeDe.Range("M3").Value = "True"
$wsDeDe.Range("N3").Value = ""
$wsDeDe.Range("O3").Value = "False"
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("A3"), $hyperlinkTarget, [Type]::Missing, [Type]::Missing, $mdName)

$loDeDe = $wsDeDe.ListObjects.Item(1)
$loDeDe.Resize($wsDeDe.Range("A1:P3"))

# ---------------------------------------------------------------------
# Column width adjustments (Status-type columns widen to fit "Ready for
# handoff"). The underlying engine quantizes ColumnWidth to steps of
# 1/6 character, so 16.33 is the closest input that reproduces the
# ~17.22 stored width seen in the target workbook.
# ---------------------------------------------------------------------
$wsOverview.Range("E1").ColumnWidth = 16.33
$wsOverview.Range("F1").ColumnWidth = 16.33
$wsZhCn.Range("C1").ColumnWidth = 16.33
$wsDeDe.Range("C1").ColumnWidth = 16.33
